$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells stay as text (avoid Excel auto-converting numeric-looking
# strings like "7.40" -> 7.4), matching the inlineStr/text cells in the source file.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.666.20'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.558.36'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.84'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.486'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.53'
$ws.Range('E8').Value = '  +3.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.245'
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0584'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.778.25'
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.557.60'
$ws.Range('E13').Value = '  -1.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.641.47'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.514'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.64'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.32'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.89'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.40'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0673'
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.91'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.98'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  +2.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.64'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.80'
$ws.Range('E26').Value = '  -1.25%  '
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.25'
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0460'
$ws.Range('E30').Value = '  -4.11%  '
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.17'
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.393.88'
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.00'
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.05'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.49'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.69'
$ws.Range('E37').Value = '  +1.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.28'
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  +3.00%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.518'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.997'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.777'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0465'
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.33'
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.31'
$ws.Range('E46').Value = '  -1.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.693.39'
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.867'
$ws.Range('E48').Value = '  -6.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.94'
$ws.Range('E49').Value = '  +6.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.35'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  -1.09%  '
